# Recruitment Test Case - OrangeHRM.xlsx
#
# The "Klik tombol Edit (ikon pensil)" positive test case (originally row 17,
# i.e. TC_RC_013) was removed. The test cases that followed it each move up
# one row (their Scenario/Steps/Test Data/Expected/Type/Status shift up,
# while the Test Case ID column stays sequential TC_RC_001..TC_RC_018), and
# the now-unused last row of that table becomes blank. The "Negative Test
# Case" table below (starting at row 25) is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 18-23's Scenario..Status (columns B-G) up into rows 17-22,
# effectively deleting the "Edit" test case row and pulling every
# subsequent positive test case up by one, without disturbing column A
# (Test Case ID, which stays sequential) or any row below the table.
for ($r = 17; $r -le 22; $r++) {
    $srcRow = $r + 1
    for ($c = 2; $c -le 7; $c++) {
        $ws.Cells.Item($r, $c).Value = $ws.Cells.Item($srcRow, $c).Value2
    }
}

# The last row of the table (previously row 23) no longer has a test case;
# clear it completely (contents + formatting) so it disappears.
$ws.Range("A23:G23").Clear()
